# Rerun the analysis scripts: add ci.lower / ci.upper confidence-interval
# columns (G, H) to the ml_results output sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-touch the full G:H block with the default style so every row in
# the existing table (including the NA/blank rows 12-14) gets a real
# (but empty) cell, matching the other sparse columns C:F.
$ws.Range("G1:H26").Style = "Normal"

# Header row
$ws.Range("G1").Value() = "ci.lower"
$ws.Range("H1").Value() = "ci.upper"

$ws.Range("G2").Value() = 0.455548133391735
$ws.Range("H2").Value() = 0.85007536142025
$ws.Range("G3").Value() = -0.0226030678004197
$ws.Range("H3").Value() = 0.0781042976342173
$ws.Range("G4").Value() = -0.124145631754953
$ws.Range("H4").Value() = -0.00998405217652044
$ws.Range("G5").Value() = -0.115945834381858
$ws.Range("H5").Value() = 0.400647736757226
$ws.Range("G6").Value() = -0.54169651842182
$ws.Range("H6").Value() = -0.0435643705486008
$ws.Range("G7").Value() = -0.106586980509125
$ws.Range("H7").Value() = 0.368308466935741
$ws.Range("G8").Value() = -0.585420888394294
$ws.Range("H8").Value() = -0.0470807761201827
$ws.Range("G9").Value() = -0.071539513934685
$ws.Range("H9").Value() = 0.032225286885847
$ws.Range("G10").Value() = 0.10942059518633
$ws.Range("H10").Value() = 0.165901590833547
$ws.Range("G11").Value() = 0.0661645636298681
$ws.Range("H11").Value() = 0.123466350135403
# row 12: estimate unavailable in source model -> leave G12/H12 blank
# row 13: estimate unavailable in source model -> leave G13/H13 blank
# row 14: estimate unavailable in source model -> leave G14/H14 blank
$ws.Range("G15").Value() = 0.0223441325578198
$ws.Range("H15").Value() = 0.127972554158613
$ws.Range("G16").Value() = 0.105365946006206
$ws.Range("H16").Value() = 0.603467115890957
$ws.Range("G17").Value() = 0.0732533014739278
$ws.Range("H17").Value() = 0.646429996574334
$ws.Range("G18").Value() = -0.0452061356008395
$ws.Range("H18").Value() = 0.156208595268435
$ws.Range("G19").Value() = -0.21317396101825
$ws.Range("H19").Value() = 0.736616933871483
$ws.Range("G20").Value() = -0.231891668763716
$ws.Range("H20").Value() = 0.801295473514452
$ws.Range("G21").Value() = 0.0661645636298681
$ws.Range("H21").Value() = 0.123466350135403
$ws.Range("G22").Value() = -0.0644505737716941
$ws.Range("H22").Value() = 0.14307902786937
$ws.Range("G23").Value() = 0.312005481569212
$ws.Range("H23").Value() = 0.582217669371882
$ws.Range("G24").Value() = -0.303922994482984
$ws.Range("H24").Value() = 0.674703172570844
$ws.Range("G25").Value() = 0.302905254330432
$ws.Range("H25").Value() = 0.567057537015356
$ws.Range("G26").Value() = -0.339689823164839
$ws.Range("H26").Value() = 0.640248809759892
